$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 7955
$ws.Range("F9").Value = 58
$ws.Range("F10").Value = 6874
$ws.Range("F13").Value = 474
$ws.Range("F15").Value = 683
$ws.Range("F16").Value = 346
$ws.Range("F22").Value = 11220
$ws.Range("F25").Value = 2134
$ws.Range("F26").Value = 2903
$ws.Range("F29").Value = 2541
$ws.Range("F34").Value = 2298
$ws.Range("F36").Value = 1563
$ws.Range("F38").Value = 69
$ws.Range("F39").Value = 5664
$ws.Range("F40").Value = 72
$ws.Range("F41").Value = 1238
$ws.Range("F42").Value = 805
$ws.Range("F46").Value = 1052
$ws.Range("F47").Value = 1480
$ws.Range("F49").Value = 1122

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 247
$ws.Range("F20").Value = 59

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 164
$ws.Range("F3").Value = 278

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 164
$ws.Range("F5").Value = 278
$ws.Range("F8").Value = 7955
$ws.Range("F11").Value = 58
$ws.Range("F12").Value = 6874
$ws.Range("F13").Value = 6874
$ws.Range("F16").Value = 474
$ws.Range("F17").Value = 683
$ws.Range("F18").Value = 346
$ws.Range("F21").Value = 247
$ws.Range("F25").Value = 11220
$ws.Range("F28").Value = 2134
$ws.Range("F29").Value = 2903
$ws.Range("F30").Value = 2541
$ws.Range("F35").Value = 2298
$ws.Range("F37").Value = 1563
$ws.Range("F39").Value = 69
$ws.Range("F40").Value = 5664
$ws.Range("F41").Value = 59
$ws.Range("F42").Value = 1238
$ws.Range("F43").Value = 805
$ws.Range("F47").Value = 1052
$ws.Range("F48").Value = 1480
$ws.Range("F50").Value = 1122
